# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, row 2
$wsOverview.Range("G2").Value = "2016-09-03 09:09:01"

# zh-cn sheet: "Correspond Handoff Datetime" column H and
# "Correspond Handback DateTime" column K, row 2
$wsZhCn.Range("H2").Value = "2016-09-03 09:08:56"
$wsZhCn.Range("K2").Value = "2016-09-03 09:09:17"

# de-de sheet: "Correspond Handback DateTime" column K, row 2
$wsDeDe.Range("K2").Value = "2016-09-03 09:09:24"
